# Parabank register update:
#   - Zipcode (F2), Phone (G2), SSN (H2) and Username (I2) on the "Demo User"
#     row become plain text fields (format "@") instead of numbers with the
#     old 00000 / phone / SSN masks.
#   - Username (I2) changes from "tcbdemotestuser2" to "tcbdemotestuser11".
#   - Column I widens a bit to fit the longer username.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# Text-format the four cells first so the values below are stored as text,
# matching the target shared-string cells rather than numeric cells.
$ws.Range("F2:I2").NumberFormat = "@"

$ws.Range("F2").Value = "123456"
$ws.Range("G2").Value = "1234567890"
$ws.Range("H2").Value = "1234567890"
$ws.Range("I2").Value = "tcbdemotestuser11"

# Column I needs to grow to fit "tcbdemotestuser11".
$ws.Columns.Item(9).ColumnWidth = 18.6
